$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: NINDYA RIZQY - update tanggal lahir dan alamat
$ws.Range("E2").Value = "Jakarta, 2025-10-29"
$ws.Range("F2").Value = "RT04/RW04, KLEYANG GUNUNG, PUNGANGAN, BAWOLATO, KABUPATEN NIAS, SUMATERA UTARA"

# Row 3: replace IDA FANIA with LIA ETIKASARI (full candidate record)
$ws.Range("A3").Value = "LIA ETIKASARI"
$ws.Range("B3").Value = "Perempuan"
$ws.Range("C3").Value = "liaetikasari0826@gmail.com"
$ws.Range("D3").Value = 6285158040206
$ws.Range("E3").Value = "TEMANGGUNG, 2007-01-11"
$ws.Range("F3").Value = "JL MELATI, BULU, KABUPATEN TEMANGGUNG, JAWA TENGAH"
$ws.Range("G3").Value = "'098776445667"
$ws.Range("H3").Value = 12345
$ws.Range("I3").Value = "AKTIF"

# Row 4: replace LIA ETIKASARI (old row4) with KURNIA AINUN (new candidate record)
$ws.Range("A4").Value = "KURNIA AINUN"
$ws.Range("B4").Value = "Perempuan"
$ws.Range("C4").Value = "kurnia@gmail.com"
$ws.Range("D4").Value = 6256667889001
$ws.Range("E4").Value = "TEMANGGUNG, 2025-10-01"
$ws.Range("F4").Value = "KRAJAN, KANDANGAN, KABUPATEN TEMANGGUNG, JAWA TENGAH"
$ws.Range("G4").Value = "'0998877765554"
$ws.Range("H4").Value = 90909
$ws.Range("I4").Value = "AKTIF"

# Re-fit column widths (Tanggal lahir / Alamat / Kode pos) for the updated content
$ws.Columns.Item(5).ColumnWidth = 25.1666666666667
$ws.Columns.Item(6).ColumnWidth = 92.1666666666667
$ws.Columns.Item(7).ColumnWidth = 15.1666666666667
